$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D and E columns as Text format before assigning, to preserve the
# original non-numeric string representations (e.g. "30.894.94", "1.002",
# "  +0.57%  ") exactly as plain text, matching the source inline strings.
$deCells = @(
"D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $deCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = '30.894.94'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '1.924.74'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '240.81'
$ws.Range("E5").Value = '  -2.67%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.4914'
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("D8").Value = '0.2982'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("D9").Value = '0.06786'
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '1.926.15'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = '17.15'
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").Value = '0.07316'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("E13").Value = '  +2.35%  '
$ws.Range("D14").Value = '89.85'
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D15").Value = '0.6752'
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").Value = '30.876.79'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '0.000008001'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '2.144.53'
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '5.218'
$ws.Range("E22").Value = '  +7.79%  '
$ws.Range("D23").Value = '205.76'
$ws.Range("E23").Value = '  +7.79%  '
$ws.Range("D24").Value = '6.309'
$ws.Range("E24").Value = '  +3.96%  '
$ws.Range("D25").Value = '9.706'
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("D26").Value = '159.40'
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").Value = '19.11'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '2.001'
$ws.Range("E28").Value = '  +4.50%  '
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").Value = '4.364'
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").Value = '4.092'
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").Value = '0.05219'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '0.7578'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").Value = '2.726'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '0.01867'
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").Value = '2.746'
$ws.Range("E38").Value = '  +2.65%  '
$ws.Range("D39").Value = '0.9293'
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").Value = '2.103'
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("D41").Value = '0.4537'
$ws.Range("E41").Value = '  +2.22%  '
$ws.Range("D42").Value = '108.41'
$ws.Range("E42").Value = '  +2.52%  '
$ws.Range("D43").Value = '5.925'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '70.34'
$ws.Range("E45").Value = '  +21.40%  '
$ws.Range("D46").Value = '0.1404'
$ws.Range("E46").Value = '  +4.48%  '
$ws.Range("D47").Value = '7.724'
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.137'
$ws.Range("E48").Value = '  +4.88%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '35.53'
$ws.Range("E49").Value = '  +5.95%  '
$ws.Range("D50").Value = '0.4112'
$ws.Range("E50").Value = '  +4.14%  '
$ws.Range("D51").Value = '0.05962'
$ws.Range("E51").Value = '  +1.99%  '
